# Auto-generated edit script: updates computed price/profit values
# across all eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 626.1875
$ws.Range("J17").Value = 621.24
$ws.Range("L17").Value = 1863.72
$ws.Range("N17").Value = -2199.72
# Row 57
$ws.Range("H57").Value = 135188.28
$ws.Range("J57").Value = 135188.28
$ws.Range("L57").Value = 405564.84
$ws.Range("N57").Value = -406562.84
# Row 92
$ws.Range("H92").Value = 197.44444
$ws.Range("I92").Value = 118.14286
$ws.Range("J92").Value = 475
$ws.Range("K92").Value = 118.14286
$ws.Range("L92").Value = 475
$ws.Range("M92").Value = 1129.85714
$ws.Range("N92").Value = -2971

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5115.317
$ws.Range("I32").Value = 3485.0535
$ws.Range("K32").Value = 3485.0535
$ws.Range("M32").Value = -3198.0535
# Row 63
$ws.Range("H63").Value = 2933.318
$ws.Range("I63").Value = 1953.3
$ws.Range("J63").Value = 3750
$ws.Range("K63").Value = 1953.3
$ws.Range("L63").Value = 3750
$ws.Range("M63").Value = -1267.3
$ws.Range("N63").Value = -5122
# Row 66
$ws.Range("H66").Value = 2933.318
$ws.Range("I66").Value = 1953.3
$ws.Range("J66").Value = 3750
$ws.Range("K66").Value = 9766.5
$ws.Range("L66").Value = 18750
$ws.Range("M66").Value = -6334.5
$ws.Range("N66").Value = -25614
# Row 102
$ws.Range("H102").Value = 4248.4
$ws.Range("J102").Value = 2997.3333
$ws.Range("L102").Value = 2997.3333
$ws.Range("N102").Value = -6241.3333
# Row 110
$ws.Range("H110").Value = 3375.6155
$ws.Range("I110").Value = 3444.0908
$ws.Range("J110").Value = 2999
$ws.Range("K110").Value = 3444.0908
$ws.Range("L110").Value = 2999
$ws.Range("M110").Value = -1399.0908
$ws.Range("N110").Value = -7089

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")
# Row 133
$ws.Range("H133").Value = 86199.664
$ws.Range("J133").Value = 85445
$ws.Range("L133").Value = 85445
$ws.Range("N133").Value = -95565
# Row 134
$ws.Range("H134").Value = 3642.9744
$ws.Range("I134").Value = 3353.4055
$ws.Range("J134").Value = 9000
$ws.Range("K134").Value = 10060.2165
$ws.Range("L134").Value = 27000
$ws.Range("M134").Value = -7525.216499999999
$ws.Range("N134").Value = -32070

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1679.5217
$ws.Range("I22").Value = 1661.8182
$ws.Range("J22").Value = 1695.75
$ws.Range("K22").Value = 1661.8182
$ws.Range("L22").Value = 1695.75
$ws.Range("M22").Value = -1311.8182
$ws.Range("N22").Value = -2395.75
# Row 28
$ws.Range("H28").Value = 6250
$ws.Range("J28").Value = 6250
$ws.Range("L28").Value = 6250
$ws.Range("N28").Value = -6740
# Row 31
$ws.Range("H31").Value = 3954.9312
$ws.Range("I31").Value = 3536.0625
$ws.Range("K31").Value = 3536.0625
$ws.Range("M31").Value = -3241.0625
# Row 34
$ws.Range("H34").Value = 3954.9312
$ws.Range("I34").Value = 3536.0625
$ws.Range("K34").Value = 3536.0625
$ws.Range("M34").Value = -3334.0625
# Row 43
$ws.Range("H43").Value = 96270.71000000001
$ws.Range("J43").Value = 96270.71000000001
$ws.Range("L43").Value = 96270.71000000001
$ws.Range("N43").Value = -96638.71000000001
# Row 101
$ws.Range("H101").Value = 96270.71000000001
$ws.Range("J101").Value = 96270.71000000001
$ws.Range("L101").Value = 96270.71000000001
$ws.Range("N101").Value = -102760.71
# Row 110
$ws.Range("H110").Value = 89999.5
$ws.Range("J110").Value = 89999.5
$ws.Range("L110").Value = 89999.5
$ws.Range("N110").Value = -98179.5
# Row 132
$ws.Range("H132").Value = 3710.125
$ws.Range("I132").Value = 3433.262
$ws.Range("K132").Value = 10299.786
$ws.Range("M132").Value = -7769.786
# Row 134
$ws.Range("H134").Value = 6808.4365
$ws.Range("I134").Value = 7001.4565
$ws.Range("K134").Value = 21004.3695
$ws.Range("M134").Value = -18469.3695

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 2495.4
$ws.Range("I129").Value = 594.4286
$ws.Range("J129").Value = 4158.75
$ws.Range("K129").Value = 1783.2858
$ws.Range("L129").Value = 12476.25
$ws.Range("M129").Value = 3216.7142
$ws.Range("N129").Value = -22476.25
# Row 131
$ws.Range("H131").Value = 1018749.2
$ws.Range("J131").Value = 8660.571
$ws.Range("L131").Value = 25981.713
$ws.Range("N131").Value = -36061.713
# Row 137
$ws.Range("H137").Value = 12008.667
$ws.Range("J137").Value = 12716.706
$ws.Range("L137").Value = 38150.118
$ws.Range("N137").Value = -48350.118

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 21.181818
$ws.Range("I2").Value = 21.181818
$ws.Range("K2").Value = 21.181818
$ws.Range("M2").Value = 91.81818200000001
# Row 80
$ws.Range("H80").Value = 15488195
$ws.Range("I80").Value = 40002256
$ws.Range("J80").Value = 5629.1577
$ws.Range("K80").Value = 40002256
$ws.Range("L80").Value = 5629.1577
$ws.Range("M80").Value = -40001258
$ws.Range("N80").Value = -7625.1577
# Row 83
$ws.Range("H83").Value = 15488195
$ws.Range("I83").Value = 40002256
$ws.Range("J83").Value = 5629.1577
$ws.Range("K83").Value = 200011280
$ws.Range("L83").Value = 28145.7885
$ws.Range("M83").Value = -200006288
$ws.Range("N83").Value = -38129.7885
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("N93").ClearContents()
# Row 126
$ws.Range("H126").Value = 5319.5
$ws.Range("I126").Value = 4632.364
$ws.Range("K126").Value = 13897.092
$ws.Range("M126").Value = -11427.092
# Row 136
$ws.Range("H136").Value = 51118.555
$ws.Range("J136").Value = 51118.555
$ws.Range("L136").Value = 153355.665
$ws.Range("N136").Value = -158455.665

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2602.1724
$ws.Range("I46").Value = 2249.5715
$ws.Range("J46").Value = 2931.2666
$ws.Range("K46").Value = 2249.5715
$ws.Range("L46").Value = 2931.2666
$ws.Range("M46").Value = -2061.5715
$ws.Range("N46").Value = -3307.2666
# Row 82
$ws.Range("H82").Value = 395.71
$ws.Range("I82").Value = 393.5258
$ws.Range("K82").Value = 393.5258
$ws.Range("M82").Value = -32.5258
# Row 85
$ws.Range("H85").Value = 395.71
$ws.Range("I85").Value = 393.5258
$ws.Range("K85").Value = 393.5258
$ws.Range("M85").Value = 854.4742

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 6872.5
$ws.Range("I122").Value = 4787.162
$ws.Range("J122").Value = 17895
$ws.Range("K122").Value = 14361.486
$ws.Range("L122").Value = 53685
$ws.Range("M122").Value = -11911.486
$ws.Range("N122").Value = -58585
# Row 132
$ws.Range("H132").Value = 2701.2
$ws.Range("I132").Value = 2719.7632
$ws.Range("K132").Value = 8159.2896
$ws.Range("M132").Value = -5629.2896
# Row 136
$ws.Range("H136").Value = 1432.75
$ws.Range("I136").Value = 1361.6
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 4084.8
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -1534.8

